$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing column F ("PackageErrorColumn"),
# pushing it to column G, and give the new column F a "BusinessKey" header.
# Insert() copies the formatting of the column to its left, so the new
# header cell keeps the same bold/underline header style as the rest of
# row 1 without any extra style definitions being created.
$ws.Columns.Item(6).Insert()
$ws.Cells.Item(1, 6).Value = "BusinessKey"

# Leave the selection on F2, matching where the cursor ends up after
# typing the new header and pressing Enter.
$ws.Range("F2").Select()
